$wb = $excel.ActiveWorkbook

# Update "Planilha1" (sheet1) - only column C changes
$ws1 = $wb.Worksheets.Item("Planilha1")

$ws1.Range("C1").Value = "R$ 18.72"
$ws1.Range("C2").Value = "R$ 72.22"
$ws1.Range("C3").Value = "R$ 28.34"
$ws1.Range("C4").Value = "R$ 21.57"
$ws1.Range("C5").Value = "R$ 20.91"
$ws1.Range("C6").Value = "R$ 61.49"
$ws1.Range("C7").Value = "R$ 71.72"
$ws1.Range("C8").Value = "R$ 36.96"
$ws1.Range("C9").Value = "R$ 76.9"
$ws1.Range("C10").Value = "R$ 66.75"

# Update "Planilha2" (sheet2) - columns A, B, C, D change
$ws2 = $wb.Worksheets.Item("Planilha2")

$ws2.Range("A1").Value = "cód 94"
$ws2.Range("B1").Value = "R$ 93.66"
$ws2.Range("C1").Value = "R$ 18.12"
$ws2.Range("D1").Value = "R$ 80.4"

$ws2.Range("A2").Value = "cód 66"
$ws2.Range("B2").Value = "R$ 95.56"
$ws2.Range("C2").Value = "R$ 84.4"
$ws2.Range("D2").Value = "R$ 55.68"

$ws2.Range("A3").Value = "cód 87"
$ws2.Range("B3").Value = "R$ 85.3"
$ws2.Range("C3").Value = "R$ 33.41"
$ws2.Range("D3").Value = "R$ 47.3"

$ws2.Range("A4").Value = "cód 34"
$ws2.Range("B4").Value = "R$ 89.22"
$ws2.Range("C4").Value = "R$ 67.24"
$ws2.Range("D4").Value = "R$ 14.08"

$ws2.Range("A5").Value = "cód 39"
$ws2.Range("B5").Value = "R$ 44.44"
$ws2.Range("C5").Value = "R$ 14.25"
$ws2.Range("D5").Value = "R$ 24.62"

$ws2.Range("A6").Value = "cód 81"
$ws2.Range("B6").Value = "R$ 21.12"
$ws2.Range("C6").Value = "R$ 82.64"
$ws2.Range("D6").Value = "R$ 90.49"

$ws2.Range("A7").Value = "cód 49"
$ws2.Range("B7").Value = "R$ 49.6"
$ws2.Range("C7").Value = "R$ 86.09"
$ws2.Range("D7").Value = "R$ 23.6"

$ws2.Range("A8").Value = "cód 79"
$ws2.Range("B8").Value = "R$ 63.48"
$ws2.Range("C8").Value = "R$ 98.14"
$ws2.Range("D8").Value = "R$ 45.39"

$ws2.Range("A9").Value = "cód 57"
$ws2.Range("B9").Value = "R$ 70.49"
$ws2.Range("C9").Value = "R$ 41.72"
$ws2.Range("D9").Value = "R$ 29.94"

$ws2.Range("A10").Value = "cód 64"
$ws2.Range("B10").Value = "R$ 88.22"
$ws2.Range("C10").Value = "R$ 60.62"
$ws2.Range("D10").Value = "R$ 55.42"
